$wb = $excel.ActiveWorkbook

# Reference the existing "data" worksheet
$dataSheet = $wb.Worksheets.Item("data")

# Add the new "metadata" worksheet, positioned right after "data"
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (row 1)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Multiple lipomas"
$metaSheet.Range("C2").Value = 19
$metaSheet.Range("E2").Value = "2017-11-05T02:37:19.852991Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:40.461789"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/19/?format=json"

# "data_version" must be stored as the literal text "1.1" (not the number 1.1).
# Compute it as a formula result on a scratch cell and paste only the resulting
# value across, which keeps it text-typed without leaving behind any unused
# number-format styles.
$scratch = $metaSheet.Range("Z1")
$scratch.Formula = '="1.1"'
$scratch.Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$scratch.Clear()

# Copy the header cell style from the "data" sheet's header row onto the new header row
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

# Copy the style used for the leading index column
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the "time_taken" column in the "data" sheet with refreshed timestamps
$dataSheet.Range("F2").Value = "2021-10-05 14:21:40.465500"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:40.465508"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:40.465511"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:40.465513"

$dataSheet.Activate()
